$ws = $excel.ActiveWorkbook.ActiveSheet

function Set-TextValue {
    param($cellRef, $text)
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue 'D2' '69.625.87'
Set-TextValue 'E2' '  -0.69%  '
Set-TextValue 'D3' '2.501.37'
Set-TextValue 'E3' '  -0.78%  '
Set-TextValue 'E4' '  +0.06%  '
Set-TextValue 'D5' '574.37'
Set-TextValue 'E5' '  -0.58%  '
Set-TextValue 'D6' '166.34'
Set-TextValue 'E6' '  -1.22%  '
Set-TextValue 'E7' '  +0.02%  '
Set-TextValue 'E8' '  -1.41%  '
Set-TextValue 'D9' '2.500.01'
Set-TextValue 'E9' '  -0.81%  '
Set-TextValue 'E10' '  -0.99%  '
Set-TextValue 'E11' '  +0.00%  '
Set-TextValue 'D12' '0.358'
Set-TextValue 'E12' '  +3.01%  '
Set-TextValue 'E13' '  +0.63%  '
Set-TextValue 'D14' '2.958.44'
Set-TextValue 'E14' '  -0.81%  '
Set-TextValue 'D15' '69.471.78'
Set-TextValue 'E15' '  -0.47%  '
Set-TextValue 'E16' '  +0.41%  '
Set-TextValue 'D17' '24.66'
Set-TextValue 'E17' '  -1.76%  '
Set-TextValue 'D18' '2.503.91'
Set-TextValue 'E18' '  -1.20%  '
Set-TextValue 'D19' '11.19'
Set-TextValue 'E19' '  -1.36%  '
Set-TextValue 'D20' '7.50'
Set-TextValue 'E20' '  -3.98%  '
Set-TextValue 'D21' '348.34'
Set-TextValue 'E22' '  -1.23%  '
Set-TextValue 'D23' '1.93'
Set-TextValue 'E23' '  -0.85%  '
Set-TextValue 'E24' '  -0.08%  '
Set-TextValue 'D25' '70.64'
Set-TextValue 'E25' '  +2.14%  '
Set-TextValue 'E26' '  -2.34%  '
Set-TextValue 'D27' '8.72'
Set-TextValue 'E27' '  -3.23%  '
Set-TextValue 'D28' '2.630.49'
Set-TextValue 'E28' '  -0.68%  '
Set-TextValue 'D29' '1.00'
Set-TextValue 'D30' '0.0₃0889'
Set-TextValue 'E30' '  -2.36%  '
Set-TextValue 'E31' '  -0.90%  '
Set-TextValue 'D32' '457.73'
Set-TextValue 'E32' '  -2.07%  '
Set-TextValue 'E33' '  -6.14%  '
Set-TextValue 'E34' '  -1.75%  '
Set-TextValue 'E35' '  -0.09%  '
Set-TextValue 'D36' '157.10'
Set-TextValue 'E36' '  +2.71%  '
Set-TextValue 'D37' '0.116'
Set-TextValue 'E37' '  -4.29%  '
Set-TextValue 'E38' '  +0.13%  '
Set-TextValue 'D39' '18.36'
Set-TextValue 'E39' '  -0.69%  '
Set-TextValue 'E40' '  +0.01%  '
Set-TextValue 'E41' '  -1.23%  '
Set-TextValue 'E42' '  -2.45%  '
Set-TextValue 'E43' '  -0.60%  '
Set-TextValue 'E44' '  -0.47%  '
Set-TextValue 'E45' '  -5.40%  '
Set-TextValue 'D46' '1.08'
Set-TextValue 'E46' '  -7.71%  '
Set-TextValue 'D47' '140.95'
Set-TextValue 'E47' '  -1.93%  '
Set-TextValue 'E48' '  -0.73%  '
Set-TextValue 'E49' '  -3.01%  '
Set-TextValue 'D50' '0.0731'
Set-TextValue 'E50' '  -0.55%  '
Set-TextValue 'E51' '  -0.82%  '
